$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (frn_adminid, gender, team, tshirt) appended after the
# existing data which ends at row 87.
$newRows = @(
    @(323, "Male",   "Accounts",   "XLarge"),
    @(253, "Female", "Marketing",  "Medium"),
    @(85,  "Female", "Marketing",  "Medium"),
    @(153, "Female", "Operations", "Large"),
    @(139, "Female", "Accounts",   "Small"),
    @(246, "Male",   "Marketing",  "Small"),
    @(156, "Female", "Marketing",  "Medium"),
    @(237, "Male",   "Accounts",   "XLarge"),
    @(274, "Male",   "Accounts",   "Large"),
    @(308, "Female", "Biz Dev",    "Large"),
    @(200, "Female", "Accounts",   "Large")
)

$startRow = 88
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

# Reset the view back to the top-left corner / default selection so the
# scrolled position and stale selection from the original file don't
# linger in the saved workbook.
$ws.Range("A1").Select()
